# Tune some hyper parameters
#
# 1) Slide master + all slide layouts: the "datetimeFigureOut" date field
#    cached text 14/6/2019 -> 16/6/2019.
# 2) Slide 1, "TextBox 4": tweak the preprocessing step wording
#    (190 -> 380, [0,0.5] -> [0,0.25]).
# 3) Slide 3, "TextBox 117": tweak the tensor initial fill value
#    (0.5 -> 0).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text on the slide master and every slide layout.
# ---------------------------------------------------------------------------
$newDate = "16/6/2019"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 1 - "TextBox 4": adjust the divide-by value + target range.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$tb1 = $slide1.Shapes.Item(2)
$tr1 = $tb1.TextFrame.TextRange
$full1 = $tr1.Text

# TextRange.Text normalizes curly quotes to straight quotes on read, so the
# search needle must use straight quotes; what we write back uses the real
# curly quotes so the underlying run text matches the original typography.
$LDQ = [char]0x201C
$RDQ = [char]0x201D
$old1 = " Convert the value for `"Timestamp`" into minutes and divide the value by 15mins as the interval between 2 timestamp is 15mins. Divide the new value by 190 so that the value is between [0,0.5]."
$new1 = " Convert the value for " + $LDQ + "Timestamp" + $RDQ + " into minutes and divide the value by 15mins as the interval between 2 timestamp is 15mins. Divide the new value by 380 so that the value is between [0,0.25]."

$pos1 = $full1.IndexOf($old1)
if ($pos1 -ge 0) {
    $range1 = $tr1.Characters($pos1 + 1, $old1.Length)
    $range1.Text = $new1
}

# ---------------------------------------------------------------------------
# 3) Slide 3 - "TextBox 117": tensor fill value 0.5 -> 0.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$tb2 = $slide3.Shapes.Item(57)
$tr2 = $tb2.TextFrame.TextRange
$old2 = "Initialize with a 25 x 5 tensor filled with 0.5 "
$new2 = "Initialize with a 25 x 5 tensor filled with 0 "
if ($tr2.Text -eq $old2) {
    $tr2.Text = $new2
}
